# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy the formatting from the existing header cell (AC1, style "1": bold,
# centered, thin border) onto the new header cells before setting their text,
# so the new headers look consistent with the rest of the header row.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-46): team record for every player row ---
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 80   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 81   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
